# Applies the "Automatic update of files" change for the HAMMARÖ workbook:
#   1. The 32 data rows (rows 2-33) are re-ordered into a new permutation
#      (the underlying data for each case stays intact - only its row
#      position changes).
#   2. Every row's "Förändrad" date (column C) is bumped from 46073 to 46074.
#
# Strategy: snapshot every source row's data (A:R via Value2 to keep full
# floating point precision for the serial-date columns, S:Z via Formula to
# keep the HYPERLINK() formulas) BEFORE writing anything back, then write
# each captured row into its new destination row according to the mapping.
# Finally, stamp column C with the new "Förändrad" value for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 33

# new row number -> old (source) row number
$mapping = @{
  2=2; 3=4; 4=3; 5=6; 6=5; 7=7; 8=8; 9=11; 10=9; 11=12; 12=10;
  13=13; 14=14; 15=33; 16=22; 17=16; 18=17; 19=31; 20=25; 21=24;
  22=32; 23=28; 24=29; 25=23; 26=15; 27=21; 28=30; 29=19; 30=20;
  31=26; 32=27; 33=18
}

# --- Step 1: snapshot every row's current contents before mutating anything ---
$snapValues = @{}
$snapFormulas = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapValues[$r] = $ws.Range("A$r`:R$r").Value2
    $snapFormulas[$r] = $ws.Range("S$r`:Z$r").Formula
}

# --- Step 2: write each row into its new position ---
foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $ws.Range("A$newRow`:R$newRow").Value2 = $snapValues[$oldRow]
    $ws.Range("S$newRow`:Z$newRow").Formula = $snapFormulas[$oldRow]
}

# --- Step 3: bump the "Förändrad" (column C) value on every data row ---
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("C$r").Value2 = 46074
}
